$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The user cleared the "unit1" value out of cell M2 (its style is kept)
# and left M2 selected. Since "unit1" was then an unused shared string,
# Excel drops it from the shared-strings table automatically on save.
$ws.Range("M2").ClearContents() | Out-Null
$ws.Range("M2").Select() | Out-Null
